$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: headers
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "id"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "data_type"
$ws.Range("F1").Value = "value_json"
$ws.Range("G1").Value = "is_active"

# ---------------------------------------------------------------------------
# JSON payloads for value_json column (row 2-5)
# ---------------------------------------------------------------------------
$json2 = @"
[
	{
		"code": "101",
		"value": "A",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "A+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "A-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "104",
		"value": "B",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "105",
		"value": "B+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "106",
		"value": "B-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "107",
		"value": "AB",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "108",
		"value": "AB+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "109",
		"value": "AB-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "110",
		"value": "O",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "111",
		"value": "O+",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "112",
		"value": "O-",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "113",
		"value": "Don't Know",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "114",
		"value": "Not Applicable",
		"langCode": "eng",
		"active": true
	}
]
"@

$json3 = @"
[
	{
		"code": "101",
		"value": "Single",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "Married",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "Widowed",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "104",
		"value": "Divorced",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "105",
		"value": "Legally Separated",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "106",
		"value": "Annulled",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "107",
		"value": "Nullified",
		"langCode": "eng",
		"active": true
	}
]
"@

$json4 = @"
[
	{
		"code": "Document-based",
		"value": "Document-based",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "Introducer-based",
		"value": "Introducer-based",
		"langCode": "eng",
		"active": true
	}
]
"@

$json5 = @"
[
	{
		"code": "101",
		"value": "Pick-up",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "102",
		"value": "Delivery to permanent address",
		"langCode": "eng",
		"active": true
	},
	{
		"code": "103",
		"value": "Delivery to present address",
		"langCode": "eng",
		"active": true
	}
]
"@

# ---------------------------------------------------------------------------
# Row 2: bloodType
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = 10001
$ws.Range("C2").Value = "bloodType"
$ws.Range("D2").Value = "Blood Type"
$ws.Range("E2").Value = "string"
$ws.Range("F2").Value = $json2
$ws.Range("G2").Value = $true

# ---------------------------------------------------------------------------
# Row 3: maritalStatus
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = 10002
$ws.Range("C3").Value = "maritalStatus"
$ws.Range("D3").Value = "Marital Status"
$ws.Range("E3").Value = "string"
$ws.Range("F3").Value = $json3
$ws.Range("G3").Value = $true

# ---------------------------------------------------------------------------
# Row 4: registrationType
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "eng"
$ws.Range("B4").Value = 10003
$ws.Range("C4").Value = "registrationType"
$ws.Range("D4").Value = "Registration Type"
$ws.Range("E4").Value = "string"
$ws.Range("F4").Value = $json4
$ws.Range("G4").Value = $true

# ---------------------------------------------------------------------------
# Row 5: modeOfClaim
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "eng"
$ws.Range("B5").Value = 10004
$ws.Range("C5").Value = "modeOfClaim"
$ws.Range("D5").Value = "Mode of Claim"
$ws.Range("E5").Value = "string"
$ws.Range("F5").Value = $json5
$ws.Range("G5").Value = $true

# ---------------------------------------------------------------------------
# Style: column A (lang_code data cells) reuses the bold/bordered header
# style, matching s="1" on A2:A5 in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---------------------------------------------------------------------------
# The multi-line value_json text otherwise triggers an implicit custom row
# height; AutoFit brings every data row back to the sheet's normal height
# (no ht/customHeight attribute), matching the target workbook. Each row is
# auto-fitted individually (a combined multi-row range does not reset the
# customHeight flag).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(4).AutoFit() | Out-Null
$ws.Rows.Item(5).AutoFit() | Out-Null

Write-Host "Done"
